$d = $word.ActiveDocument

$replacements = @(
    @{ old = "39×13="; new = "45×20=" },
    @{ old = "63×90="; new = "47×14=" },
    @{ old = "91×25="; new = "83×78=" },
    @{ old = "89×69="; new = "61×45=" },
    @{ old = "87×69="; new = "38×45=" },
    @{ old = "72×88="; new = "46×64=" },
    @{ old = "85×70="; new = "26×88=" },
    @{ old = "98×52="; new = "90×70=" },
    @{ old = "18×11="; new = "31×69=" },
    @{ old = "61×65="; new = "68×26=" },
    @{ old = "62×59="; new = "87×93=" },
    @{ old = "52×88="; new = "77×11=" },
    @{ old = "48×91="; new = "46×72=" },
    @{ old = "48×13="; new = "87×82=" },
    @{ old = "99×89="; new = "49×35=" },
    @{ old = "90×22="; new = "89×18=" },
    @{ old = "65×29="; new = "21×96=" },
    @{ old = "51×66="; new = "11×14=" },
    @{ old = "73×74="; new = "95×61=" },
    @{ old = "52×16="; new = "96×62=" },
    @{ old = "73×33="; new = "40×92=" },
    @{ old = "90×86="; new = "88×58=" },
    @{ old = "66×64="; new = "11×14=" },
    @{ old = "25×26="; new = "94×53=" },
    @{ old = "39×31="; new = "25×97=" }
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $r.new, 2)
}
